# Refresh the crypto price/volume table (rows 2-51) with the latest
# scrape results. Mirrors the "Updated cryptos list ... with GitHub
# Actions" automation commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some refreshed prices are plain decimal-looking strings (e.g. "1.035").
# Format those Price cells as Text first so Excel keeps them as literal
# strings instead of silently converting them to numbers.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.741.27"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("D3").Value = "1.907.62"
$ws.Range("E3").Value = "  +2.70%  "
$ws.Range("D4").Value = "1.035"
$ws.Range("E4").Value = "  +3.05%  "
$ws.Range("D5").Value = "320.55"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("D6").Value = "1.031"
$ws.Range("E6").Value = "  +2.72%  "
$ws.Range("D7").Value = "0.5202"
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("D8").Value = "0.3970"
$ws.Range("E8").Value = "  +3.63%  "
$ws.Range("D9").Value = "0.08396"
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("D10").Value = "1.137"
$ws.Range("E10").Value = "  +2.61%  "
$ws.Range("D11").Value = "42.59"
$ws.Range("E11").Value = "  +2.59%  "
$ws.Range("D12").Value = "6.306"
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").Value = "1.912.38"
$ws.Range("E13").Value = "  +2.72%  "
$ws.Range("D14").Value = "20.66"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").Value = "7.322"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").Value = "1.033"
$ws.Range("E16").Value = "  +2.94%  "
$ws.Range("D17").Value = "0.00001113"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "91.59"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "0.06801"
$ws.Range("E19").Value = "  +2.33%  "
$ws.Range("D20").Value = "18.02"
$ws.Range("E20").Value = "  +2.01%  "
$ws.Range("D21").Value = "1.032"
$ws.Range("E21").Value = "  +2.91%  "
$ws.Range("D22").Value = "6.115"
$ws.Range("E22").Value = "  +1.99%  "
$ws.Range("D23").Value = "28.801.68"
$ws.Range("E23").Value = "  +2.83%  "
$ws.Range("D24").Value = "11.28"
$ws.Range("E24").Value = "  +2.07%  "
$ws.Range("D25").Value = "2.295"
$ws.Range("E25").Value = "  +2.30%  "
$ws.Range("D26").Value = "2.131.14"
$ws.Range("E26").Value = "  +2.80%  "
$ws.Range("D27").Value = "163.42"
$ws.Range("E27").Value = "  +3.45%  "
$ws.Range("D28").Value = "21.05"
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("D29").Value = "2.480"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").Value = "127.93"
$ws.Range("E30").Value = "  +2.85%  "
$ws.Range("D31").Value = "0.1073"
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("D32").Value = "1.054"
$ws.Range("E32").Value = "  +2.52%  "
$ws.Range("D33").Value = "5.962"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D34").Value = "3.681"
$ws.Range("E34").Value = "  +2.28%  "
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "0.02478"
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "9.460"
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("D37").Value = "0.06640"
$ws.Range("E37").Value = "  +2.41%  "
$ws.Range("D38").Value = "0.2235"
$ws.Range("E38").Value = "  +3.05%  "
$ws.Range("D39").Value = "0.6580"
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").Value = "1.262"
$ws.Range("E40").Value = "  +3.61%  "
$ws.Range("D41").Value = "1.200"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").Value = "5.030"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").Value = "0.6186"
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D45").Value = "13.27"
$ws.Range("E45").Value = "  +2.21%  "
$ws.Range("D46").Value = "3.767"
$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("D47").Value = "1.306"
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("D48").Value = "2.023"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("D49").Value = "1.242"
$ws.Range("E49").Value = "  +2.28%  "
$ws.Range("D50").Value = "122.76"
$ws.Range("E50").Value = "  +2.10%  "
$ws.Range("D51").Value = "0.06980"
$ws.Range("E51").Value = "  +2.16%  "
